$d = $word.ActiveDocument

# --- 1. Swap the two embedded image names (image3.jpg <-> image1.jpg) ---
$shape1 = $d.Shapes.Item(1)
$shape2 = $d.Shapes.Item(2)
$shape1.Name = "image1.jpg"
$shape2.Name = "image3.jpg"

# --- 2. Append a sentence to the end of the "content management" paragraph ---
$p43 = $d.Paragraphs.Item(43)
$p43End = $p43.Range.End
$insPoint = $d.Range($p43End - 1, $p43End - 1)
$insPoint.InsertAfter(" really feel like these should’ve been included in the lesson material…")

# --- 3. Insert two new paragraphs after it: one blank, one with the SEO text ---
$p43 = $d.Paragraphs.Item(43)
$p43.Range.InsertParagraphAfter()

$p44 = $d.Paragraphs.Item(44)
$p44.Range.InsertParagraphAfter()

$p45 = $d.Paragraphs.Item(45)
$r45 = $p45.Range
$r45.Font.NameAscii = "Helvetica Neue Light"
$r45.Font.NameFarEast = "Helvetica Neue Light"
$r45.Font.NameBi = "Helvetica Neue Light"
$r45.Font.NameOther = "Helvetica Neue Light"
$r45.Font.Size = 14
$r45.Font.SizeBi = 14
$r45.InsertAfter("In terms of SEO, I did my best to ensure that all pages have unique descriptions to the extent possible to my knowledge. I also made sure the site tags are specific for this project and tried to include key search terms in posts where possible and reasonable to do so.")

# --- 4. Fill in the "What was difficult" blank paragraph ---
$p47 = $d.Paragraphs.Item(47)
$r47 = $p47.Range
$r47.Font.NameAscii = "Helvetica Neue Light"
$r47.Font.NameFarEast = "Helvetica Neue Light"
$r47.Font.NameBi = "Helvetica Neue Light"
$r47.Font.NameOther = "Helvetica Neue Light"
$r47.Font.Size = 14
$r47.Font.SizeBi = 14
$r47.InsertAfter("Content Management wasn’t too easy to begin with. As I mentioned above, I encountered a number of issues with content loading when trying to use vanilla WordPress. I’ve gotten very used to exams being quite strict in public universities here in Norway, and thus felt like going outside of the tools provided almost feels a bit like cheating. Because of that, I didn’t really look into plugins until Connor O’Brien mentioned CPT and ACF. I know this is a silly mindset to have and something I’m working to get rid of, but it’s not always immediately obvious what is considered “cheating” and what isn’t when it comes to these projects. I realize Bootstrap at this stage is frowned upon, but that’s about where the obvious rules stop at least in my mind. I would’ve liked a more detailed run-down of how to go about solving these issues other than “Tutor support will be limited”.")

# --- 5. Fill in the "What would you do differently" blank paragraph ---
$p49 = $d.Paragraphs.Item(49)
$r49 = $p49.Range
$r49.Font.NameAscii = "Helvetica Neue Light"
$r49.Font.NameFarEast = "Helvetica Neue Light"
$r49.Font.NameBi = "Helvetica Neue Light"
$r49.Font.NameOther = "Helvetica Neue Light"
$r49.Font.Size = 14
$r49.Font.SizeBi = 14
$r49.InsertAfter("Next time, I will do more research on possible tools to use. I want to get into JavaScript Libraries, but didn’t find any which would’ve been good to use with this project. I’ve previously been given feedback on certain topics to look into following a project, which I tried to do this time around. I still feel a bit like I don’t quite know how or when to use these kinds of tools, but it’s slowly getting there. I do feel like I have made a project I can be proud of though, and definitely a project which has challenged me and my abilities - and in the end made me a better developer.")
